$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H3").Value = 0
$ws.Range("H4").Value = 2
$ws.Range("H5").Value = 0
$ws.Range("H6").Value = 0
$ws.Range("H7").Value = 2
$ws.Range("H8").Value = 0
$ws.Range("H9").Value = 0
$ws.Range("H10").Value = 0
$ws.Range("H13").Value = 0
$ws.Range("H14").Value = 2
$ws.Range("H15").Value = 2
$ws.Range("H17").Value = 2
$ws.Range("H18").Value = 2
$ws.Range("H22").Value = 2
$ws.Range("H23").Value = 0
$ws.Range("H24").Value = 2
$ws.Range("H25").Value = 0
$ws.Range("H27").Value = 0
$ws.Range("H28").Value = 2
$ws.Range("H29").Value = 2
$ws.Range("H31").Value = 0
$ws.Range("H32").Value = 0
$ws.Range("H33").Value = 0
$ws.Range("H35").Value = 2
$ws.Range("H37").Value = 2
$ws.Range("H38").Value = 2
$ws.Range("H39").Value = 0
$ws.Range("H40").Value = 2
$ws.Range("H41").Value = 2
$ws.Range("H42").Value = 2
$ws.Range("H43").Value = 0
$ws.Range("H44").Value = 0
$ws.Range("H45").Value = 0
$ws.Range("H46").Value = 2
$ws.Range("H47").Value = 2
$ws.Range("H48").Value = 0
$ws.Range("H50").Value = 2
$ws.Range("H51").Value = 0
$ws.Range("H53").Value = 2
$ws.Range("H55").Value = 0
$ws.Range("H56").Value = 0
$ws.Range("H58").Value = 2
$ws.Range("H59").Value = 0
$ws.Range("H60").Value = 2
$ws.Range("H61").Value = 2
$ws.Range("H62").Value = 2
$ws.Range("H63").Value = 1.05
